$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4188.846
$ws.Range("I28").Value = 910
$ws.Range("K28").Value = 910
$ws.Range("M28").Value = -425
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H98").Value = 890.0476
$ws.Range("I98").Value = 890.0476
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 890.0476
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 607.9524
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 1783.7693
$ws.Range("I100").Value = 1118.5
$ws.Range("K100").Value = 1118.5
$ws.Range("M100").Value = -577.5
$ws.Range("H122").Value = 890.0476
$ws.Range("I122").Value = 890.0476
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2670.1428
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -220.1428000000001
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 52635108
$ws.Range("I137").Value = 35717200
$ws.Range("J137").Value = 100005256
$ws.Range("K137").Value = 107151600
$ws.Range("L137").Value = 300015768
$ws.Range("M137").Value = -107149050
$ws.Range("N137").Value = -300020868
$ws.Range("H138").Value = 8066.756
$ws.Range("J138").Value = 9864.134
$ws.Range("L138").Value = 29592.402
$ws.Range("N138").Value = -39872.402

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1538487.9
$ws.Range("I11").Value = 4000002.5
$ws.Range("J11").Value = 41.25
$ws.Range("K11").Value = 4000002.5
$ws.Range("L11").Value = 41.25
$ws.Range("M11").Value = -3999858.5
$ws.Range("N11").Value = -329.25
$ws.Range("H32").Value = 17866558
$ws.Range("I32").Value = 21282388
$ws.Range("K32").Value = 21282388
$ws.Range("M32").Value = -21282101
$ws.Range("H61").Value = 71434504
$ws.Range("I61").Value = 125003096
$ws.Range("K61").Value = 125003096
$ws.Range("M61").Value = -125002884
$ws.Range("H74").Value = 333708380
$ws.Range("I74").Value = 333708380
$ws.Range("K74").Value = 333708380
$ws.Range("M74").Value = -333707506
$ws.Range("H77").Value = 333708380
$ws.Range("I77").Value = 333708380
$ws.Range("K77").Value = 1668541900
$ws.Range("M77").Value = -1668537532
$ws.Range("H102").Value = 2851.6155
$ws.Range("I102").Value = 2964.3333
$ws.Range("K102").Value = 2964.3333
$ws.Range("M102").Value = -1342.3333
$ws.Range("H110").Value = 15969.77
$ws.Range("I110").Value = 15969.77
$ws.Range("K110").Value = 15969.77
$ws.Range("M110").Value = -13924.77
$ws.Range("H132").Value = 90914040
$ws.Range("I132").Value = 5442.2
$ws.Range("J132").Value = 1000000000
$ws.Range("K132").Value = 16326.6
$ws.Range("L132").Value = 3000000000
$ws.Range("M132").Value = -13796.6
$ws.Range("N132").Value = -3000005060
$ws.Range("H136").Value = 71434504
$ws.Range("I136").Value = 125003096
$ws.Range("K136").Value = 375009288
$ws.Range("M136").Value = -375006738

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3849.1936
$ws.Range("I99").Value = 2501.9443
$ws.Range("K99").Value = 2501.9443
$ws.Range("M99").Value = -1003.9443
$ws.Range("H134").Value = 1615.6595
$ws.Range("I134").Value = 1615.6595
$ws.Range("K134").Value = 4846.9785
$ws.Range("M134").Value = -2311.9785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20837940
$ws.Range("I31").Value = 3430.3225
$ws.Range("K31").Value = 3430.3225
$ws.Range("M31").Value = -3135.3225
$ws.Range("H34").Value = 20837940
$ws.Range("I34").Value = 3430.3225
$ws.Range("K34").Value = 3430.3225
$ws.Range("M34").Value = -3228.3225
$ws.Range("H86").Value = 6149.3335
$ws.Range("I86").Value = 4907
$ws.Range("J86").Value = 6397.8
$ws.Range("K86").Value = 4907
$ws.Range("L86").Value = 6397.8
$ws.Range("M86").Value = -3784
$ws.Range("N86").Value = -8643.799999999999
$ws.Range("H89").Value = 6149.3335
$ws.Range("I89").Value = 4907
$ws.Range("J89").Value = 6397.8
$ws.Range("K89").Value = 24535
$ws.Range("L89").Value = 31989
$ws.Range("M89").Value = -18919
$ws.Range("N89").Value = -43221
$ws.Range("H109").Value = 61428.145
$ws.Range("J109").Value = 62999.418
$ws.Range("L109").Value = 62999.418
$ws.Range("N109").Value = -65079.418
$ws.Range("H141").Value = 299072.47
$ws.Range("J141").Value = 307599.72
$ws.Range("L141").Value = 307599.72
$ws.Range("N141").Value = -317959.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I2").Value = 1390.8334
$ws.Range("J2").Value = 3806036.5
$ws.Range("K2").Value = 8345.000400000001
$ws.Range("L2").Value = 22836219
$ws.Range("M2").Value = -8232.000400000001
$ws.Range("N2").Value = -22836445
$ws.Range("H17").Value = 109.75
$ws.Range("J17").Value = 99.5
$ws.Range("L17").Value = 298.5
$ws.Range("N17").Value = -636.5
$ws.Range("H56").Value = 12067
$ws.Range("I56").Value = 12067
$ws.Range("K56").Value = 12067
$ws.Range("M56").Value = -11537
$ws.Range("H75").Value = 2722.077
$ws.Range("I75").Value = 2248.6667
$ws.Range("J75").Value = 3127.8572
$ws.Range("K75").Value = 6746.000100000001
$ws.Range("L75").Value = 9383.571599999999
$ws.Range("M75").Value = -5748.000100000001
$ws.Range("N75").Value = -11379.5716
$ws.Range("H78").Value = 2722.077
$ws.Range("I78").Value = 2248.6667
$ws.Range("J78").Value = 3127.8572
$ws.Range("K78").Value = 20238.0003
$ws.Range("L78").Value = 28150.7148
$ws.Range("M78").Value = -15246.0003
$ws.Range("N78").Value = -38134.7148
$ws.Range("H122").Value = 2279.2144
$ws.Range("I122").Value = 1713.625
$ws.Range("J122").Value = 3033.3333
$ws.Range("K122").Value = 15422.625
$ws.Range("L122").Value = 27299.9997
$ws.Range("M122").Value = -12972.625
$ws.Range("N122").Value = -32199.9997
$ws.Range("H134").Value = 4232.684
$ws.Range("J134").Value = 19499.334
$ws.Range("L134").Value = 58498.00199999999
$ws.Range("N134").Value = -68638.00199999999
$ws.Range("H137").Value = 4513.875
$ws.Range("I137").Value = 2572.3
$ws.Range("K137").Value = 7716.900000000001
$ws.Range("M137").Value = -2616.900000000001
$ws.Range("H139").Value = 3398.1052
$ws.Range("I139").Value = 2833.2
$ws.Range("K139").Value = 8499.599999999999
$ws.Range("M139").Value = -3359.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1038.4546
$ws.Range("J107").Value = 699.6667
$ws.Range("L107").Value = 699.6667
$ws.Range("N107").Value = -4539.6667
$ws.Range("H132").Value = 3822.6086
$ws.Range("I132").Value = 3899.5625
$ws.Range("K132").Value = 11698.6875
$ws.Range("M132").Value = -9168.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H93").Value = 1706
$ws.Range("I93").Value = 1407
$ws.Range("K93").Value = 1407
$ws.Range("M93").Value = -159
$ws.Range("H100").Value = 3460.4285
$ws.Range("I100").Value = 3274.5715
$ws.Range("K100").Value = 3274.5715
$ws.Range("M100").Value = -2733.5715
$ws.Range("H136").Value = 5942.9
$ws.Range("I136").Value = 5298.3887
$ws.Range("J136").Value = 11743.5
$ws.Range("K136").Value = 15895.1661
$ws.Range("L136").Value = 35230.5
$ws.Range("M136").Value = -13345.1661
$ws.Range("N136").Value = -40330.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 391.125
$ws.Range("I113").Value = 242.15384
$ws.Range("K113").Value = 726.4615200000001
$ws.Range("M113").Value = 1443.53848
$ws.Range("H132").Value = 5070.357
$ws.Range("I132").Value = 4836.694
$ws.Range("J132").Value = 6706
$ws.Range("K132").Value = 14510.082
$ws.Range("L132").Value = 20118
$ws.Range("M132").Value = -11980.082
$ws.Range("N132").Value = -25178
$ws.Range("H136").Value = 2437.125
$ws.Range("I136").Value = 1584
$ws.Range("J136").Value = 4996.5
$ws.Range("K136").Value = 4752
$ws.Range("L136").Value = 14989.5
$ws.Range("M136").Value = -2202
$ws.Range("N136").Value = -20089.5
